$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 57
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 11
$ws.Range("H4").Value = 11
$ws.Range("E5").Value = 62
$ws.Range("E6").Value = 21
$ws.Range("E7").Value = 13
$ws.Range("E10").Value = 183
$ws.Range("F10").Value = 75
$ws.Range("H10").Value = 75
$ws.Range("E11").Value = 134
$ws.Range("F11").Value = 65
$ws.Range("H11").Value = 65
$ws.Range("E12").Value = 194
$ws.Range("F12").Value = 98
$ws.Range("H12").Value = 98
$ws.Range("E13").Value = 64
$ws.Range("F13").Value = 27
$ws.Range("H13").Value = 27
$ws.Range("E14").Value = 57
$ws.Range("E15").Value = 83
$ws.Range("E16").Value = 76
$ws.Range("F16").Value = 33
$ws.Range("H16").Value = 33
$ws.Range("E17").Value = 36
$ws.Range("F17").Value = 18
$ws.Range("H17").Value = 18
$ws.Range("E18").Value = 28
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = 13
$ws.Range("E21").Value = 64
$ws.Range("E22").Value = 78
$ws.Range("E23").Value = 86
$ws.Range("E24").Value = 90
$ws.Range("F24").Value = 42
$ws.Range("H24").Value = 42
$ws.Range("E25").Value = 79
$ws.Range("E26").Value = 47
$ws.Range("F26").Value = 20
$ws.Range("H26").Value = 20
$ws.Range("E27").Value = 127
$ws.Range("F27").Value = 59
$ws.Range("H27").Value = 59
$ws.Range("E28").Value = 82
$ws.Range("E30").Value = 87
$ws.Range("E31").Value = 36
$ws.Range("E32").Value = 87
$ws.Range("F32").Value = 43
$ws.Range("H32").Value = 43
$ws.Range("E33").Value = 116
$ws.Range("E34").Value = 93
$ws.Range("F34").Value = 48
$ws.Range("H34").Value = 48
$ws.Range("E35").Value = 58
$ws.Range("F35").Value = 29
$ws.Range("H35").Value = 29
$ws.Range("E36").Value = 30
$ws.Range("E37").Value = 58
$ws.Range("E38").Value = 42
$ws.Range("F38").Value = 23
$ws.Range("H38").Value = 23
$ws.Range("E39").Value = 95
$ws.Range("E40").Value = 119
$ws.Range("F40").Value = 47
$ws.Range("H40").Value = 47
$ws.Range("E41").Value = 161
$ws.Range("F41").Value = 55
$ws.Range("H41").Value = 55
$ws.Range("E42").Value = 137
$ws.Range("E43").Value = 40
$ws.Range("E44").Value = 126
$ws.Range("E46").Value = 104
$ws.Range("F46").Value = 45
$ws.Range("H46").Value = 45
$ws.Range("E47").Value = 180
$ws.Range("F47").Value = 70
$ws.Range("H47").Value = 70
$ws.Range("E48").Value = 91
$ws.Range("E49").Value = 100
$ws.Range("E50").Value = 81
$ws.Range("E51").Value = 82
$ws.Range("E52").Value = 4
